$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new mapped field "gpm_id" under the "Tabela Nome novo" column (B),
# continuing the list of renamed fields for the distribuicao_horas_extras table.
$ws.Range("B10").Value = "gpm_id"

# Move/save the active selection to B11, matching the cursor position left
# after entering the new value.
$ws.Range("B11").Select()
